$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot with the latest values.
# Columns D/E are plain text cells; force text format first when the new
# value would otherwise be auto-parsed as a number by Excel.
$ws.Range("D2").Value = '43.086.21'
$ws.Range("E2").Value = '  +4.81%  '
$ws.Range("D3").Value = '2.239.51'
$ws.Range("E3").Value = '  +2.99%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.64'
$ws.Range("E5").Value = '  +3.61%  '
$ws.Range("E6").Value = '  +0.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.73'
$ws.Range("E7").Value = '  +7.81%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  +6.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.13'
$ws.Range("E10").Value = '  +2.83%  '
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.48'
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("E13").Value = '  +3.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.101'
$ws.Range("E14").Value = '  +0.13%  '
$ws.Range("D15").Value = '2.555.45'
$ws.Range("E15").Value = '  +2.24%  '
$ws.Range("E16").Value = '  +5.58%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.245.06'
$ws.Range("E17").Value = '  +4.15%  '
$ws.Range("B18").Value = 'Polygon'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.815'
$ws.Range("E18").Value = '  +0.58%  '
$ws.Range("D19").Value = '42.992.08'
$ws.Range("E19").Value = '  +4.93%  '
$ws.Range("E20").Value = '  +3.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.17'
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.00'
$ws.Range("E22").Value = '  +0.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.54'
$ws.Range("E23").Value = '  +4.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.24'
$ws.Range("E24").Value = '  +15.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '229.06'
$ws.Range("E25").Value = '  +1.40%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E28").Value = '  -4.65%  '
$ws.Range("E29").Value = '  +1.96%  '
$ws.Range("E30").Value = '  +1.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.81'
$ws.Range("E31").Value = '  +22.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.86'
$ws.Range("E32").Value = '  +4.15%  '
$ws.Range("E33").Value = '  +1.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0795'
$ws.Range("E34").Value = '  +3.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.40'
$ws.Range("E35").Value = '  +4.48%  '
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("E37").Value = '  +7.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.38'
$ws.Range("E38").Value = '  +6.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0332'
$ws.Range("E39").Value = '  +16.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '13.25'
$ws.Range("E40").Value = '  +6.77%  '
$ws.Range("E41").Value = '  +3.36%  '
$ws.Range("E42").Value = '  +2.97%  '
$ws.Range("E43").Value = '  +4.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.28'
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '105.43'
$ws.Range("E45").Value = '  +7.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.59'
$ws.Range("E46").Value = '  +3.05%  '
$ws.Range("E47").Value = '  +2.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.445'
$ws.Range("E48").Value = '  +19.37%  '
$ws.Range("E49").Value = '  +1.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.32'
$ws.Range("E50").Value = '  +3.73%  '
$ws.Range("E51").Value = '  +1.49%  '
